$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(0.3301136626125567,0.3949069889952284,0.2832973020472173,0.3571804451842914,0.3017493060930475,7,7,7,7,7,2.310795638287897,2.764348922966599,1.983081114330521,2.50026311629004,2.112245142651332,63.88048,80.64758792200021,48.10386154325745,72.54372444754591,55.00654197501315,0.203616653773568,0.2687802782889044,0.1494746064757829,0.2357895216797287,0.1723380961865035,0.8595319604176113,0.9218073099881111,0.7929035260783918,0.8933913491789289,0.822971758626735)

$row3 = @(0.33,0.33,0.33,0.33,0.33,7.008650454676773,9.137338641110436,5.443269109730695,7.917152192314232,6.061317022253602,2.312854650043335,3.015321751566444,1.796278806211129,2.612660223463696,2.000234617343688,63.8488,75.21552451439547,55.42784932428966,69.18625840385873,58.90981425312337,0.202436494678118,0.300552560485084,0.116192452246413,0.2524199123774446,0.1526893179694813,0.8528999012808594,0.9440024691536069,0.7341154569326024,0.9052864392390564,0.7914206999631612)

$row4 = @(0.3300789700090945,0.3943883010179898,0.2833683820696578,0.3570496897887267,0.3017349428211397,7.006449683098982,9.12798705085791,5.431402308235739,7.911817098952892,6.059075348741236,2.312727900809143,3.155615094726688,1.703826992705189,2.666498669303904,1.943636146755307,65.10521,90.68672292392961,47.36474784479361,75.00734373692225,54.9017405427993,0.2013971555282028,0.3169361038270316,0.09917400691363822,0.2602388622497128,0.142669912107736,0.8479147069084788,0.9544819082535033,0.6971593999335332,0.9097051757792216,0.7726048598037696)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value2 = $row2[$i]
}

for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value2 = $row3[$i]
}

for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value2 = $row4[$i]
}
